# Update weekly epidemiological poisson data (semana 50 de 2025)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - evento 113
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 2

# Row 4 - evento 115
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0.37

# Row 5 - evento 155
$ws.Range("C5").Value = 8
$ws.Range("E5").Value = 0.14

# Row 6 - evento 210
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 0.15

# Row 9 - evento 300
$ws.Range("C9").Value = 48
$ws.Range("D9").Value = 49

# Row 10 - evento 330
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0.37

# Row 12 - evento 342
$ws.Range("C12").Value = 4
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 0.07000000000000001

# Row 13 - evento 346
$ws.Range("C13").Value = 11

# Row 14 - evento 348
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0.27

# Row 17 - evento 356
$ws.Range("C17").Value = 8
$ws.Range("E17").Value = 0.12

# Row 19 - evento 365
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 0.01

# Row 22 - evento 455
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0.37

# Row 23 - evento 465
$ws.Range("D23").Value = 2

# Row 25 - evento 549
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 0

# Row 26 - evento 560
$ws.Range("C26").Value = 2
$ws.Range("E26").Value = 0.14

# Row 29 - evento 620
$ws.Range("C29").Value = 1
$ws.Range("E29").Value = 0.37

# Row 32 - evento 750
$ws.Range("C32").Value = 1
$ws.Range("D32").Value = 2
$ws.Range("E32").Value = 0.18

# Row 35 - evento 831
$ws.Range("C35").Value = 10
$ws.Range("D35").Value = 1
$ws.Range("E35").Value = 0

# Row 36 - evento 850
$ws.Range("C36").Value = 10
$ws.Range("D36").Value = 13
$ws.Range("E36").Value = 0.07000000000000001

# Row 37 - new event 895 Zika
# Use a leading apostrophe so Excel stores "895" as text (not a number),
# then reset the style so no extra formatting is attached to the cell.
$ws.Range("A37").Value = "'895"
$ws.Range("A37").Style = "Normal"
$ws.Range("B37").Value = "Zika"
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 1
$ws.Range("E37").Value = 0
